$d = $word.ActiveDocument

# Locate the specific "Skills" table that contains "Ada, C++" (several
# similar project-summary tables exist in this resume, one per project).
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    try {
        $txt = $candidate.Cell(7, 2).Range.Text
    } catch {
        $txt = ""
    }
    if ($txt.StartsWith("Ada, C++")) {
        $table = $candidate
        break
    }
}

# --- Edit 1: "Ada, C++" -> split into two runs "C" and ", C++" -----------
$langCell = $table.Cell(7, 2)
$langStart = $langCell.Range.Start

# Replace "Ada" (the first 3 characters) with "C".
$adaRange = $d.Range($langStart, $langStart + 3)
$adaRange.Text = "C"

# The cell text is now "C, C++". Force a run boundary between "C" and
# ", C++" by toggling a character formatting property on and back off on
# just the ", C++" portion, which keeps the visible formatting unchanged
# but prevents the engine from coalescing the two pieces into one run.
$commaRange = $d.Range($langStart + 1, $langStart + 1 + 5)
$commaRange.Bold = 1
$commaRange.Bold = 0

# --- Edit 2: "ObjectADA" + ", g++" -> single run "g++" --------------------
$toolsCell = $table.Cell(9, 2)
$toolsStart = $toolsCell.Range.Start

# Remove the leading "ObjectADA, " (11 characters), leaving just "g++".
$removeRange = $d.Range($toolsStart, $toolsStart + 11)
$removeRange.Text = ""
